$word.UserName = "Kathryn Pasqualucci"
$d = $word.ActiveDocument
$d.TrackRevisions = $true
$r = $d.Range($d.Content.End, $d.Content.End)
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:rPr><w:ins w:id="100" w:author="Kathryn Pasqualucci" w:date="2015-08-12T09:35:00Z"/></w:rPr></w:pPr><w:ins w:id="101" w:author="Kathryn Pasqualucci" w:date="2015-08-12T09:35:00Z"><w:r><w:t xml:space="preserve"> </w:t></w:r></w:ins></w:p>
<w:p><w:pPr><w:rPr><w:ins w:id="102" w:author="Kathryn Pasqualucci" w:date="2015-08-12T09:35:00Z"/></w:rPr></w:pPr></w:p>
<w:p><w:pPr><w:rPr><w:ins w:id="103" w:author="Kathryn Pasqualucci" w:date="2015-08-12T09:35:00Z"/></w:rPr></w:pPr></w:p>
<w:p><w:ins w:id="104" w:author="Kathryn Pasqualucci" w:date="2015-08-12T09:35:00Z"><w:r><w:t>I DON'T LIKE SPIDERS</w:t></w:r></w:ins></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xml)

$lastP = $d.Paragraphs.Item($d.Paragraphs.Count)
Write-Output ("lastP range " + $lastP.Range.Start + "," + $lastP.Range.End)
Write-Output ("lastP text=[" + $lastP.Range.Text + "]")

$fixXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:ins w:id="104" w:author="Kathryn Pasqualucci" w:date="2015-08-12T09:35:00Z"><w:r><w:t>I DON'T LIKE SPIDERS</w:t></w:r></w:ins></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$lastP.Range.InsertXML($fixXml)

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
